# Case_2_133 (380 kV) vm_pu results update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.038911043077652
$ws.Range("D2").Value = 1.042404171706192
$ws.Range("E2").Value = 1.054705471781406
$ws.Range("F2").Value = 1.060965707906929
$ws.Range("I2").Value = 1.041008804662115
$ws.Range("J2").Value = 1.044005915439567
$ws.Range("K2").Value = 1.045180955100124
$ws.Range("L2").Value = 1.05744794576024
$ws.Range("M2").Value = 1.063691052343577
$ws.Range("N2").Value = 1.01856014259076

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.03994571168146
$ws.Range("D3").Value = 1.043208193177583
$ws.Range("E3").Value = 1.056020900551643
$ws.Range("F3").Value = 1.062366259014263
$ws.Range("I3").Value = 1.041322895781242
$ws.Range("J3").Value = 1.044685145701347
$ws.Range("K3").Value = 1.045795819041985
$ws.Range("L3").Value = 1.058575388386204
$ws.Range("M3").Value = 1.064904654572703
$ws.Range("N3").Value = 1.018788306446662

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.04061470169102
$ws.Range("D4").Value = 1.043727901231037
$ws.Range("E4").Value = 1.056872493778473
$ws.Range("F4").Value = 1.063272942260616
$ws.Range("I4").Value = 1.041524482446008
$ws.Range("J4").Value = 1.045123563728434
$ws.Range("K4").Value = 1.046192478676446
$ws.Range("L4").Value = 1.059304795314647
$ws.Range("M4").Value = 1.065689838205654
$ws.Range("N4").Value = 1.018935504936624

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.040895823899982
$ws.Range("D5").Value = 1.043946255441363
$ws.Range("E5").Value = 1.057230607336025
$ws.Range("F5").Value = 1.063654218242449
$ws.Range("I5").Value = 1.041608834313958
$ws.Range("J5").Value = 1.045307614287109
$ws.Range("K5").Value = 1.046358947713282
$ws.Range("L5").Value = 1.059611410573143
$ws.Range("M5").Value = 1.066019908196712
$ws.Range("N5").Value = 1.018997282216682

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.040943018473727
$ws.Range("D6").Value = 1.043982910390022
$ws.Range("E6").Value = 1.05729074233663
$ws.Range("F6").Value = 1.06371824247148
$ws.Range("I6").Value = 1.041622974195533
$ws.Range("J6").Value = 1.045338501929029
$ws.Range("K6").Value = 1.046386881793219
$ws.Range("L6").Value = 1.059662891081598
$ws.Range("M6").Value = 1.066075327223564
$ws.Range("N6").Value = 1.019007648740935

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.04061845853488
$ws.Range("D7").Value = 1.043730819405604
$ws.Range("E7").Value = 1.056877278498457
$ws.Range("F7").Value = 1.063278036467731
$ws.Range("I7").Value = 1.041525611112341
$ws.Range("J7").Value = 1.045126024043612
$ws.Range("K7").Value = 1.046194704170011
$ws.Range("L7").Value = 1.059308892427663
$ws.Range("M7").Value = 1.065694248697498
$ws.Range("N7").Value = 1.018936330820008

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.039260820238153
$ws.Range("D8").Value = 1.0426760079391
$ws.Range("E8").Value = 1.055149941063048
$ws.Range("F8").Value = 1.061438942660979
$ws.Range("I8").Value = 1.041115295232846
$ws.Range("J8").Value = 1.044235690301386
$ws.Range("K8").Value = 1.045388999540449
$ws.Range("L8").Value = 1.057828997260195
$ws.Range("M8").Value = 1.064101216602942
$ws.Range("N8").Value = 1.018637342535583

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.036864545309502
$ws.Range("D9").Value = 1.04081309861732
$ws.Range("E9").Value = 1.052109276945612
$ws.Range("F9").Value = 1.058201420968949
$ws.Range("I9").Value = 1.040379611823242
$ws.Range("J9").Value = 1.042658445853765
$ws.Range("K9").Value = 1.04396005497664
$ws.Range("L9").Value = 1.055220184424684
$ws.Range("M9").Value = 1.06129323603498
$ws.Range("N9").Value = 1.018107122607234

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.035264320691502
$ws.Range("D10").Value = 1.039568330199702
$ws.Range("E10").Value = 1.050084108500756
$ws.Range("F10").Value = 1.056045052878393
$ws.Range("I10").Value = 1.039880635844456
$ws.Range("J10").Value = 1.041601297284982
$ws.Range("K10").Value = 1.043001227810107
$ws.Range("L10").Value = 1.053480129579491
$ws.Range("M10").Value = 1.059420529078927
$ws.Range("N10").Value = 1.017751374320942

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.034570749789636
$ws.Range("D11").Value = 1.039028658119541
$ws.Range("E11").Value = 1.049207612408062
$ws.Range("F11").Value = 1.055111752951329
$ws.Range("I11").Value = 1.039662549511996
$ws.Range("J11").Value = 1.041142192201742
$ws.Range("K11").Value = 1.042584570129253
$ws.Range("L11").Value = 1.052726433944384
$ws.Range("M11").Value = 1.058609422601158
$ws.Range("N11").Value = 1.017596791653954

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.034313025759595
$ws.Range("D12").Value = 1.038828097502085
$ws.Range("E12").Value = 1.048882101924778
$ws.Range("F12").Value = 1.054765143777901
$ws.Range("I12").Value = 1.03958123772208
$ws.Range("J12").Value = 1.040971456148142
$ws.Range("K12").Value = 1.042429582149555
$ws.Range("L12").Value = 1.052446438883318
$ws.Range("M12").Value = 1.058308106768867
$ws.Range("N12").Value = 1.017539291256664

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.034368313026433
$ws.Range("D13").Value = 1.038871123073796
$ws.Range("E13").Value = 1.048951922401132
$ws.Range("F13").Value = 1.054839489962539
$ws.Range("I13").Value = 1.039598693176507
$ws.Range("J13").Value = 1.04100808885395
$ws.Range("K13").Value = 1.042462837686563
$ws.Range("L13").Value = 1.052506500582096
$ws.Range("M13").Value = 1.058372741658186
$ws.Range("N13").Value = 1.017551628978672

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.034549448307084
$ws.Range("D14").Value = 1.039012081806842
$ws.Range("E14").Value = 1.049180704420942
$ws.Range("F14").Value = 1.055083100925028
$ws.Range("I14").Value = 1.039655834474355
$ws.Range("J14").Value = 1.041128083271832
$ws.Range("K14").Value = 1.042571763321959
$ws.Range("L14").Value = 1.052703290276094
$ws.Range("M14").Value = 1.058584516466711
$ws.Range("N14").Value = 1.01759204031708

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.034661038320525
$ws.Range("D15").Value = 1.039098917557962
$ws.Range("E15").Value = 1.049321672322091
$ws.Range("F15").Value = 1.055233205525244
$ws.Range("I15").Value = 1.039691000711809
$ws.Range("J15").Value = 1.04120198874076
$ws.Range("K15").Value = 1.042638846466279
$ws.Range("L15").Value = 1.052824533656408
$ws.Range("M15").Value = 1.058714993220637
$ws.Range("N15").Value = 1.017616928269687

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.035310336627259
$ws.Range("D16").Value = 1.039604132120346
$ws.Range("E16").Value = 1.050142287135646
$ws.Range("F16").Value = 1.056107001516976
$ws.Range("I16").Value = 1.039895066790161
$ws.Range("J16").Value = 1.041631738012467
$ws.Range("K16").Value = 1.043028848797912
$ws.Range("L16").Value = 1.053530144570399
$ws.Range("M16").Value = 1.059474354828648
$ws.Range("N16").Value = 1.017761622042801

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.035717445863195
$ws.Range("D17").Value = 1.039920857802078
$ws.Range("E17").Value = 1.050657146042095
$ws.Range("F17").Value = 1.056655221172035
$ws.Range("I17").Value = 1.040022529257054
$ws.Range("J17").Value = 1.041900945676749
$ws.Range("K17").Value = 1.043273090561203
$ws.Range("L17").Value = 1.053972689298759
$ws.Range("M17").Value = 1.059950623483745
$ws.Range("N17").Value = 1.017852239568906

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.035954841902236
$ws.Range("D18").Value = 1.040105533030014
$ws.Range("E18").Value = 1.050957494807235
$ws.Range("F18").Value = 1.05697502933527
$ws.Range("I18").Value = 1.040096680460088
$ws.Range("J18").Value = 1.04205783947966
$ws.Range("K18").Value = 1.043415409965244
$ws.Range("L18").Value = 1.054230795145763
$ws.Range("M18").Value = 1.060228402672919
$ws.Range("N18").Value = 1.017905043022653

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.036035777019028
$ws.Range("D19").Value = 1.040168491410101
$ws.Range("E19").Value = 1.051059912902821
$ws.Range("F19").Value = 1.057084082696687
$ws.Range("I19").Value = 1.040121930955407
$ws.Range("J19").Value = 1.042111314104014
$ws.Range("K19").Value = 1.043463913018566
$ws.Range("L19").Value = 1.054318798753256
$ws.Range("M19").Value = 1.060323114850352
$ws.Range("N19").Value = 1.017923038784688

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.035673773512637
$ws.Range("D20").Value = 1.039886882911561
$ws.Range("E20").Value = 1.050601902410776
$ws.Range("F20").Value = 1.056596398183484
$ws.Range("I20").Value = 1.0400088739644
$ws.Range("J20").Value = 1.041872075742769
$ws.Range("K20").Value = 1.04324690049368
$ws.Range("L20").Value = 1.053925210837293
$ws.Range("M20").Value = 1.0598995264807
$ws.Range("N20").Value = 1.017842522568195

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.034496111277817
$ws.Range("D21").Value = 1.038970575802208
$ws.Range("E21").Value = 1.049113332210008
$ws.Range("F21").Value = 1.055011361939991
$ws.Range("I21").Value = 1.039639016203711
$ws.Range("J21").Value = 1.041092753543517
$ws.Range("K21").Value = 1.042539693590967
$ws.Range("L21").Value = 1.052645341713821
$ws.Range("M21").Value = 1.058522155047301
$ws.Range("N21").Value = 1.017580142436779

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.033755083775325
$ws.Range("D22").Value = 1.038393864650086
$ws.Range("E22").Value = 1.048177749925057
$ws.Range("F22").Value = 1.054015132061217
$ws.Range("I22").Value = 1.039404708196248
$ws.Range("J22").Value = 1.040601582183668
$ws.Range("K22").Value = 1.042093755421613
$ws.Range("L22").Value = 1.051840410478646
$ws.Range("M22").Value = 1.057655944330071
$ws.Range("N22").Value = 1.017414701937519

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.034147972314775
$ws.Range("D23").Value = 1.038699646423605
$ws.Range("E23").Value = 1.048673688550943
$ws.Range("F23").Value = 1.054543220767422
$ws.Range("I23").Value = 1.039529086610953
$ws.Range("J23").Value = 1.040862073568876
$ws.Range("K23").Value = 1.042330278056443
$ws.Range("L23").Value = 1.052267142191559
$ws.Range("M23").Value = 1.058115159204251
$ws.Range("N23").Value = 1.01750244985054

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.035693507365003
$ws.Range("D24").Value = 1.039902234904337
$ws.Range("E24").Value = 1.050626864504616
$ws.Range("F24").Value = 1.056622977632763
$ws.Range("I24").Value = 1.040015044807258
$ws.Range("J24").Value = 1.041885121226106
$ws.Range("K24").Value = 1.043258735098494
$ws.Range("L24").Value = 1.05394666438073
$ws.Range("M24").Value = 1.059922615078799
$ws.Range("N24").Value = 1.01784691342381

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.037484512910459
$ws.Range("D25").Value = 1.04129520337569
$ws.Range("E25").Value = 1.052895007645463
$ws.Range("F25").Value = 1.059038036558031
$ws.Range("I25").Value = 1.040571304478557
$ws.Range("J25").Value = 1.043067195399251
$ws.Range("K25").Value = 1.044330562354232
$ws.Range("L25").Value = 1.05589476387239
$ws.Range("M25").Value = 1.062019282539331
$ws.Range("N25").Value = 1.018244596505922
